$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.436.60"
$ws.Range("E2").Value = "  +0.42%  "
$ws.Range("D3").Value = "3.680.97"
$ws.Range("E3").Value = "  +0.06%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'685.91"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.92%  "
$ws.Range("D6").Value = "'159.82"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.17%  "
$ws.Range("E7").Value = "  -0.16%  "
$ws.Range("E8").Value = "  +0.13%  "
$ws.Range("E9").Value = "  -1.04%  "
$ws.Range("D10").Value = "'7.06"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.97%  "
$ws.Range("E11").Value = "  -2.98%  "
$ws.Range("E12").Value = "  -0.84%  "
$ws.Range("D13").Value = "4.303.72"
$ws.Range("E13").Value = "  +0.01%  "
$ws.Range("D14").Value = "'32.25"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.22%  "
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "3.685.61"
$ws.Range("E15").Value = "  +0.35%  "
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "69.427.57"
$ws.Range("E16").Value = "  +0.30%  "
$ws.Range("E17").Value = "  +1.98%  "
$ws.Range("D18").Value = "'15.82"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.61%  "
$ws.Range("E19").Value = "  -3.06%  "
$ws.Range("D20").Value = "'470.74"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.60%  "
$ws.Range("D21").Value = "'9.93"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.70%  "
$ws.Range("E22").Value = "  -1.47%  "
$ws.Range("D23").Value = "'79.59"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.23%  "
$ws.Range("D24").Value = "3.829.20"
$ws.Range("E24").Value = "  -0.08%  "
$ws.Range("E25").Value = "  +0.12%  "
$ws.Range("D26").Value = "'0.0000124"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.70%  "
$ws.Range("D27").Value = "'10.97"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.61%  "
$ws.Range("D28").Value = "'9.19"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.96%  "
$ws.Range("E29").Value = "  -0.79%  "
$ws.Range("E30").Value = "  -4.28%  "
$ws.Range("E31").Value = "  -4.76%  "
$ws.Range("D32").Value = "'6.56"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.95%  "
$ws.Range("D33").Value = "'1.00"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.16%  "
$ws.Range("D34").Value = "'26.86"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.50%  "
$ws.Range("D35").Value = "3.657.02"
$ws.Range("E35").Value = "  +0.23%  "
$ws.Range("E36").Value = "  -3.28%  "
$ws.Range("D37").Value = "'8.19"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.29%  "
$ws.Range("D38").Value = "'6.11"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.82%  "
$ws.Range("E39").Value = "  +0.00%  "
$ws.Range("D40").Value = "'2.21"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.97%  "
$ws.Range("D41").Value = "'0.0896"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.82%  "
$ws.Range("E43").Value = "  -1.30%  "
$ws.Range("D44").Value = "'165.39"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.34%  "
$ws.Range("D45").Value = "'47.49"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.95%  "
$ws.Range("B46").Value = "FLOKI"
$ws.Range("C46").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D46").Value = "'0.000283"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.47%  "
$ws.Range("B47").Value = "dogwifhat"
$ws.Range("C47").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D47").Value = "'2.73"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.53%  "
$ws.Range("E48").Value = "  +5.74%  "
$ws.Range("E49").Value = "  +1.30%  "
$ws.Range("D50").Value = "'27.65"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.49%  "
$ws.Range("E51").Value = "  -3.10%  "
